$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 458.44446
$ws.Range("I19").Value = 230
$ws.Range("J19").Value = 487
$ws.Range("K19").Value = 230
$ws.Range("L19").Value = 487
$ws.Range("M19").Value = -55
$ws.Range("N19").Value = -837

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3038.55
$ws.Range("J62").Value = 3351.9092
$ws.Range("L62").Value = 3351.9092
$ws.Range("N62").Value = -4599.9092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3038.55
$ws.Range("J65").Value = 3351.9092
$ws.Range("L65").Value = 16759.546
$ws.Range("N65").Value = -22999.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 736.5
$ws.Range("I2").Value = 779.1905
$ws.Range("J2").Value = 437.66666
$ws.Range("K2").Value = 779.1905
$ws.Range("L2").Value = 437.66666
$ws.Range("M2").Value = -666.1905
$ws.Range("N2").Value = -663.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3520.375
$ws.Range("I45").Value = 3599.6667
$ws.Range("K45").Value = 3599.6667
$ws.Range("M45").Value = -3222.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 736.5
$ws.Range("I116").Value = 779.1905
$ws.Range("J116").Value = 437.66666
$ws.Range("K116").Value = 779.1905
$ws.Range("L116").Value = 437.66666
$ws.Range("M116").Value = 1514.8095
$ws.Range("N116").Value = -5025.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2180.3333
$ws.Range("I122").Value = 2180.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6540.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4090.999899999999
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 22069.654
$ws.Range("I132").Value = 2738.3125
$ws.Range("J132").Value = 52999.8
$ws.Range("K132").Value = 8214.9375
$ws.Range("L132").Value = 158999.4
$ws.Range("M132").Value = -5684.9375
$ws.Range("N132").Value = -164059.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 736.5
$ws.Range("I3").Value = 779.1905
$ws.Range("J3").Value = 437.66666
$ws.Range("K3").Value = 779.1905
$ws.Range("L3").Value = 437.66666
$ws.Range("M3").Value = -665.1905
$ws.Range("N3").Value = -665.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 24386.8
$ws.Range("J35").Value = 24386.8
$ws.Range("L35").Value = 24386.8
$ws.Range("N35").Value = -25006.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4288.654
$ws.Range("I31").Value = 927.5
$ws.Range("J31").Value = 5782.5
$ws.Range("K31").Value = 927.5
$ws.Range("L31").Value = 5782.5
$ws.Range("M31").Value = -632.5
$ws.Range("N31").Value = -6372.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4288.654
$ws.Range("I34").Value = 927.5
$ws.Range("J34").Value = 5782.5
$ws.Range("K34").Value = 927.5
$ws.Range("L34").Value = 5782.5
$ws.Range("M34").Value = -725.5
$ws.Range("N34").Value = -6186.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2926.9678
$ws.Range("I99").Value = 2488.0908
$ws.Range("J99").Value = 3999.7778
$ws.Range("K99").Value = 2488.0908
$ws.Range("L99").Value = 3999.7778
$ws.Range("M99").Value = -990.0907999999999
$ws.Range("N99").Value = -6995.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1163.4706
$ws.Range("I105").Value = 973.4286
$ws.Range("J105").Value = 2050.3333
$ws.Range("K105").Value = 973.4286
$ws.Range("L105").Value = 2050.3333
$ws.Range("M105").Value = 773.5714
$ws.Range("N105").Value = -5544.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1037
$ws.Range("I107").Value = 359.4375
$ws.Range("J107").Value = 1940.4166
$ws.Range("K107").Value = 359.4375
$ws.Range("L107").Value = 1940.4166
$ws.Range("M107").Value = 1560.5625
$ws.Range("N107").Value = -5780.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2926.9678
$ws.Range("I126").Value = 2488.0908
$ws.Range("J126").Value = 3999.7778
$ws.Range("K126").Value = 7464.2724
$ws.Range("L126").Value = 11999.3334
$ws.Range("M126").Value = -4994.2724
$ws.Range("N126").Value = -16939.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 971.69446
$ws.Range("J5").Value = 1664.9
$ws.Range("L5").Value = 4994.700000000001
$ws.Range("N5").Value = -5218.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 541.2222
$ws.Range("I14").Value = 541.2222
$ws.Range("K14").Value = 1623.6666
$ws.Range("M14").Value = -1450.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 198
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 198
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 594
$ws.Range("M50").Value = $null
$ws.Range("N50").Value = -1556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 198
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 198
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 594
$ws.Range("M53").Value = $null
$ws.Range("N53").Value = -1556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1327.7273
$ws.Range("J75").Value = 1461.5
$ws.Range("L75").Value = 4384.5
$ws.Range("N75").Value = -6380.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1327.7273
$ws.Range("J78").Value = 1461.5
$ws.Range("L78").Value = 13153.5
$ws.Range("N78").Value = -23137.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5845.8335
$ws.Range("J81").Value = 5845.8335
$ws.Range("L81").Value = 17537.5005
$ws.Range("N81").Value = -19783.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 5845.8335
$ws.Range("J84").Value = 5845.8335
$ws.Range("L84").Value = 52612.5015
$ws.Range("N84").Value = -63844.5015

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 587
$ws.Range("J97").Value = 587
$ws.Range("L97").Value = 1761
$ws.Range("N97").Value = -2753

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2599.75
$ws.Range("I114").Value = 3879.6
$ws.Range("J114").Value = 466.66666
$ws.Range("K114").Value = 11638.8
$ws.Range("L114").Value = 1399.99998
$ws.Range("M114").Value = -8384.799999999999
$ws.Range("N114").Value = -7907.999980000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1086
$ws.Range("I117").Value = 821.8
$ws.Range("J117").Value = 1218.1
$ws.Range("K117").Value = 2465.4
$ws.Range("L117").Value = 3654.3
$ws.Range("M117").Value = 976.6000000000004
$ws.Range("N117").Value = -10538.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1113
$ws.Range("I122").Value = 328.6
$ws.Range("J122").Value = 1291.2727
$ws.Range("K122").Value = 2957.4
$ws.Range("L122").Value = 11621.4543
$ws.Range("M122").Value = -507.4000000000001
$ws.Range("N122").Value = -16521.4543

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2360.8
$ws.Range("I123").Value = 1145.7142
$ws.Range("J123").Value = 5196
$ws.Range("K123").Value = 3437.1426
$ws.Range("L123").Value = 15588
$ws.Range("M123").Value = -987.1425999999997
$ws.Range("N123").Value = -20488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 718.5599999999999
$ws.Range("J131").Value = 724.9583
$ws.Range("L131").Value = 2174.8749
$ws.Range("N131").Value = -12254.8749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 771.3333
$ws.Range("I132").Value = 771.3333
$ws.Range("K132").Value = 6941.9997
$ws.Range("M132").Value = -4411.9997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 971.69446
$ws.Range("J135").Value = 1664.9
$ws.Range("L135").Value = 14984.1
$ws.Range("N135").Value = -20054.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 5343.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4456.9697
$ws.Range("I113").Value = 5765.619
$ws.Range("J113").Value = 2166.8333
$ws.Range("K113").Value = 5765.619
$ws.Range("L113").Value = 2166.8333
$ws.Range("M113").Value = -3595.619
$ws.Range("N113").Value = -6506.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2599.5112
$ws.Range("I126").Value = 2158.7273
$ws.Range("J126").Value = 3021.1304
$ws.Range("K126").Value = 6476.1819
$ws.Range("L126").Value = 9063.3912
$ws.Range("M126").Value = -4006.1819
$ws.Range("N126").Value = -14003.3912

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1823.1578
$ws.Range("I100").Value = 1203.909
$ws.Range("K100").Value = 1203.909
$ws.Range("M100").Value = -662.9090000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 894306.9399999999
$ws.Range("I122").Value = 1034844.9
$ws.Range("J122").Value = 4233.3335
$ws.Range("K122").Value = 3104534.7
$ws.Range("L122").Value = 12700.0005
$ws.Range("M122").Value = -3102084.7
$ws.Range("N122").Value = -17600.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 19750
$ws.Range("J130").Value = 19750
$ws.Range("L130").Value = 19750
$ws.Range("N130").Value = -29790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1392.04
$ws.Range("I136").Value = 1392.04
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4176.12
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1626.12
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1284.9131
$ws.Range("I113").Value = 1432.4736
$ws.Range("J113").Value = 584
$ws.Range("K113").Value = 4297.4208
$ws.Range("L113").Value = 1752
$ws.Range("M113").Value = -2127.4208
$ws.Range("N113").Value = -6092
